$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all": insert a new data row (row 30) before the trailing footnote
# row, shifting the footnote down to row 31.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate()

[void]$wsAll.Rows("30:30").Insert()
$wsAll.Range("A30").Value = 43958
$wsAll.Range("B30").Value = 272
$wsAll.Range("C30").Value = 268
$wsAll.Range("D30").Value = 101
$wsAll.Range("E30").Value = 92
$wsAll.Range("F30").Value = 9
$wsAll.Range("G30").Value = 8
$wsAll.Range("H30").Value = 159

[void]$wsAll.Range("G41").Select()

# ---------------------------------------------------------------------------
# Sheet "kobe": insert a new data row (row 85) before the trailing footnote
# row, shifting the footnote down to row 86.
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Activate()

[void]$wsKobe.Rows("85:85").Insert()
$wsKobe.Range("A85").Value = 43958
$wsKobe.Range("B85").Value = 0
$wsKobe.Range("C85").Value = 2351
$wsKobe.Range("D85").Value = 0
$wsKobe.Range("E85").Value = 272
$wsKobe.Range("F85").Value = 96
$wsKobe.Range("G85").Value = 88
$wsKobe.Range("H85").Value = 8
$wsKobe.Range("I85").Value = 8
$wsKobe.Range("J85").Value = 152

[void]$wsKobe.Range("D85").Select()

# ---------------------------------------------------------------------------
# Sheet "other": insert a new data row (row 60) before the trailing footnote
# row, shifting the footnote down to row 61.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate()

[void]$wsOther.Rows("60:60").Insert()
$wsOther.Range("A60").Value = 43958
$wsOther.Range("B60").Value = 0
$wsOther.Range("C60").Value = 12
$wsOther.Range("D60").Value = 5
$wsOther.Range("E60").Value = 4
$wsOther.Range("F60").Value = 1
$wsOther.Range("G60").Value = 0
$wsOther.Range("H60").Value = 7
# The row above (59) carries a formatted-but-empty I cell that Insert()
# copies down onto I60; the new row has no data out to column I, so drop
# that inherited ghost cell entirely (ClearContents keeps the styled-empty
# cell around, Clear() removes it).
[void]$wsOther.Range("I60").Clear()

[void]$wsOther.Range("B60:H60").Select()

$wsAll.Activate()
